$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (C/D/E columns) per "rene fine" contingency re-run ---

# row 8
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# row 9
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# row 10
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# row 11
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# row 12
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

# row 13
$ws.Range("D13").Value = 8

# row 14
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Append two new contingency rows (16: line7, 17: line8) ---
# Copy formatting from the last existing data row (15) so the new rows
# match the table's look (bold/centered/bordered A column, etc.)
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
